# B6-PowerPoint.pptx commit replay
#
# The source commit swaps the table style applied to every table that was
# using the old default "Table_0" style ({82C1A7AC-8416-42DB-80E3-ABDC0E86A762})
# for a different built-in table style ({3EEABC39-5C1F-4B02-A572-9495F4439E78}).
# There are exactly three such tables in the deck (on the "comparing two
# companies" slides). We walk every slide/shape, find any table still on the
# old style id, and re-apply the new style id via Table.ApplyStyle, which is
# the supported, non-throwing way to change a table's style through the
# PowerPoint object model (Table.Style is read-only / cannot be assigned a
# style id directly).

$OldStyleId = "{82C1A7AC-8416-42DB-80E3-ABDC0E86A762}"
$NewStyleId = "{3EEABC39-5C1F-4B02-A572-9495F4439E78}"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if ($shape.HasTable) {
            $tbl = $shape.Table

            if ($tbl.Style -eq $OldStyleId) {
                $tbl.ApplyStyle($NewStyleId)
            }
        }
    }
}
